# Added 'Concentrated quench buffer' to the storage_medium list (Closes #49)

$wb = $excel.ActiveWorkbook

# 1) Insert the new "Concentrated quench buffer" entry into the storage_medium
#    lookup sheet, right after "Tris-EDTA" (row 11) and before "Unknown"
#    (previously row 12). Inserting a row shifts all subsequent rows down by one.
$wsStorage = $wb.Worksheets.Item("storage_medium")
$wsStorage.Rows.Item(12).Insert()
$wsStorage.Cells.Item(12, 1).Value = "Concentrated quench buffer"
$wsStorage.Cells.Item(12, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000391"

# 2) The storage_medium list now has 23 entries instead of 22, so the data
#    validation on the "storage_medium" column (M) of the main sheet must be
#    extended to cover the new row.
$wsMain = $wb.Worksheets.Item("Sample Suspension")
$rng = $wsMain.Range("M2:M1001")
$rng.Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$23"

# 3) Update the recorded modification timestamp in the .metadata sheet.
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Cells.Item(2, 3).Value = "2024-10-02T11:08:45-07:00"
